$d = $word.ActiveDocument

# Edit 1: "Hace en el año 2040" -> "Hace en el año 2030"
# (represented in the target OOXML as two runs: "Hace en el año 203" + "0")
$d.Content.Find.Execute("Hace en el año 2040", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hace en el año 2030", 2)

# Edit 2: merge " " + "subway" + "; " runs (with proofErr wrapping "subway")
# into a single run's text " subway; " (proofErr removed).
$d.Content.Find.Execute(" subway; ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " subway; ", 2)
